# Fixed tests for fund_units
# Updates the sample valuations workbook: new fund/portfolio-company test
# data, removes the stale hyperlinks that pointed at per-row instrument
# records, and re-applies an (inactive) AutoFilter database definition
# left behind by the refreshed export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Drop the old hyperlinks (the new export has none).
# ---------------------------------------------------------------------
if ($ws.Hyperlinks.Count -gt 0) {
    $ws.Hyperlinks.Delete()
}

# ---------------------------------------------------------------------
# 2. Column B (Valuation *) is blank for every data row - clear any
#    leftover (empty but styled) cells so they disappear entirely.
# ---------------------------------------------------------------------
$ws.Range("B2:B7").Clear()

# ---------------------------------------------------------------------
# 3. Re-key the Instrument (E) and Portfolio Company (D) columns with
#    the new bulk-test fixture data. Instrument is written before
#    Company so the shared-string table order matches the new export.
# ---------------------------------------------------------------------
$ws.Range("E2").Value = "Equity"
$ws.Range("E3").Value = "CCPS"
$ws.Range("E4").Value = "CCPS"
$ws.Range("E5").Value = "Equity"
$ws.Range("E6").Value = "CCPS"
$ws.Range("E7").Value = "CCPS"

$ws.Range("D2").Value = "TSTF1 Port Co 1"
$ws.Range("D3").Value = "TSTF1 Port Co 1"
$ws.Range("D4").Value = "TSTF1 Port Co 2"
$ws.Range("D5").Value = "TSTF1 Port Co 1"
$ws.Range("D6").Value = "TSTF1 Port Co 1"
$ws.Range("D7").Value = "TSTF1 Port Co 2"

# Re-assert the plain (non-hyperlink) Arial 10pt look for the data
# cells now that the hyperlink style is gone.
$dataRng = $ws.Range("C2:E7")
$dataRng.Style = "Normal"
$dataRng.Font.Name = "Arial"
$dataRng.Font.Size = 10

# ---------------------------------------------------------------------
# 4. Valuation Date (A) + Per Share Value (C) updated for the new
#    fixture rows (two valuation dates: 31 Mar 2024 and 31 Mar 2025).
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "3/31/2024"
$ws.Range("A3").Value = "3/31/2024"
$ws.Range("A4").Value = "3/31/2024"
$ws.Range("A5").Value = "3/31/2025"
$ws.Range("A6").Value = "3/31/2025"
$ws.Range("A7").Value = "3/31/2025"

$dateRng = $ws.Range("A2:A7")
$dateRng.NumberFormat = "mm-dd-yy"
$dateRng.Font.Name = "Arial"
$dateRng.Font.Size = 10

$ws.Range("C2").Value = 150
$ws.Range("C3").Value = 200
$ws.Range("C4").Value = 200
$ws.Range("C5").Value = 200
$ws.Range("C6").Value = 250
$ws.Range("C7").Value = 250

# ---------------------------------------------------------------------
# 5. The refreshed export carries a lot of blank, pre-formatted columns
#    (J:L) left over from the bulk test template - recreate that so the
#    used range matches (dimension grows out to the template's columns).
# ---------------------------------------------------------------------
$tailRng = $ws.Range("J4:L7")
$tailRng.Font.Name = "Arial"
$tailRng.Font.Size = 10
$tailRng.Font.Color = 0

# ---------------------------------------------------------------------
# 6. AutoFilter database left on the sheet (filter UI itself is off).
# ---------------------------------------------------------------------
$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$O`$7")
$fdb.Visible = $false

# ---------------------------------------------------------------------
# 7. View state: new machine / zoom level, selection parked on D1.
# ---------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 113
$ws.Range("D1").Select()
